$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.700.28"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "2.097.40"
$ws.Range("E3").Value = "  -1.28%  "

$ws.Range("E4").Value = "  -0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.53%  "

$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5129"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4400"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.13"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "

$ws.Range("D13").Value = "2.098.70"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.745"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.164"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001150"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.74%  "

$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06640"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.174"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("D23").Value = "29.746.78"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("E24").Value = "  -1.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.40%  "

$ws.Range("D26").Value = "2.349.46"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.509"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.129"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1045"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.629"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.151"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.964"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.031"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02564"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06689"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6845"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.294"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6687"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.299"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.610"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.217"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000335"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.163"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
